$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -12.89630000000001
$ws.Range("D11").Value = -8.609499999999993
$ws.Range("C12").Value = -10.9195
$ws.Range("C15").Value = -13.8242
$ws.Range("D23").Value = -8.244999999999999
$ws.Range("C27").Value = -12.4975
$ws.Range("C28").Value = -12.9987
$ws.Range("D28").Value = -8.193400000000002
$ws.Range("C31").Value = -13.52339999999999
$ws.Range("C32").Value = -13.41450000000001
$ws.Range("D32").Value = -8.2471
$ws.Range("D34").Value = -7.6819
$ws.Range("C36").Value = -12.85180000000001
$ws.Range("D36").Value = -8.877899999999997
$ws.Range("D37").Value = -8.322400000000002
$ws.Range("C38").Value = -11.88970000000001
$ws.Range("D42").Value = -8.874999999999998
$ws.Range("C46").Value = -14.73229999999999
$ws.Range("D49").Value = -8.080999999999996
$ws.Range("C54").Value = -12.5337
$ws.Range("D54").Value = -8.086300000000003
$ws.Range("C55").Value = -14.14890000000001
$ws.Range("C56").Value = -13.05259999999998
$ws.Range("C67").Value = -11.6966
$ws.Range("C69").Value = -11.4454
$ws.Range("C72").Value = -12.174
$ws.Range("C73").Value = -11.44500000000001
$ws.Range("D78").Value = -8.078900000000001
$ws.Range("D80").Value = -8.194900000000001
$ws.Range("C83").Value = -13.43140000000001
$ws.Range("C86").Value = -14.38979999999999
$ws.Range("C91").Value = -12.17740000000001
$ws.Range("C93").Value = -10.4922
$ws.Range("D97").Value = -8.551099999999996
$ws.Range("C99").Value = -12.81080000000001
$ws.Range("D99").Value = -8.293899999999999
$ws.Range("D100").Value = -8.054200000000002
$ws.Range("D101").Value = -7.624199999999997
$ws.Range("C104").Value = -12.94760000000001
$ws.Range("C105").Value = -12.70890000000001
